# Update "想去人数" (column F) counts on three sheets to reflect the
# latest generated output (gh-pages update at 456a3b4).

$wb = $excel.ActiveWorkbook

# ---- Sheet "展览" (Exhibition) ----
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value  = 74
$ws1.Range("F6").Value  = 844
$ws1.Range("F7").Value  = 421
$ws1.Range("F8").Value  = 4707
$ws1.Range("F9").Value  = 4707
$ws1.Range("F12").Value = 160
$ws1.Range("F15").Value = 117
$ws1.Range("F16").Value = 7516
$ws1.Range("F17").Value = 249
$ws1.Range("F18").Value = 128
$ws1.Range("F19").Value = 296
$ws1.Range("F22").Value = 1374
$ws1.Range("F26").Value = 19
$ws1.Range("F28").Value = 6180
$ws1.Range("F30").Value = 23
$ws1.Range("F33").Value = 449
$ws1.Range("F34").Value = 6433
$ws1.Range("F37").Value = 99
$ws1.Range("F40").Value = 14
$ws1.Range("F44").Value = 1116
$ws1.Range("F46").Value = 440
$ws1.Range("F47").Value = 2149
$ws1.Range("F48").Value = 46

# ---- Sheet "演出" (Performance) ----
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F13").Value = 146

# ---- Sheet "全部类型" (All Types) ----
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F6").Value  = 74
$ws4.Range("F8").Value  = 421
$ws4.Range("F9").Value  = 4707
$ws4.Range("F10").Value = 4707
$ws4.Range("F13").Value = 160
$ws4.Range("F16").Value = 117
$ws4.Range("F17").Value = 7516
$ws4.Range("F18").Value = 249
$ws4.Range("F19").Value = 128
$ws4.Range("F21").Value = 1374
$ws4.Range("F29").Value = 6180
$ws4.Range("F32").Value = 23
$ws4.Range("F35").Value = 449
$ws4.Range("F36").Value = 6433
$ws4.Range("F39").Value = 99
$ws4.Range("F44").Value = 1116
$ws4.Range("F46").Value = 440
$ws4.Range("F47").Value = 146
$ws4.Range("F48").Value = 2149
$ws4.Range("F49").Value = 46
